# Adds a new "Test" run row (row 3) for Logistic Regression with scaling,
# and extends the (Mean CV) Train / Test Accuracy columns (H, I) with a
# "0.000" number format down through the rest of the table (rows 3-30).
# Also widens column H and moves the active selection, matching the
# author's edit ("added test run for Logistic Regression with scaling").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tracker")

# --- Row 3: new test run -------------------------------------------------
$ws.Range("A3").Value = 44971.897916666669      # Timestamp
$ws.Range("B3").Value = "Test"                  # Modeling Type
$ws.Range("D3").Value = "LogReg"                # Model
$ws.Range("E3").Value = "scaling"                # data_enigneering
$ws.Range("H3").Value = 0.709                   # (Mean CV) Train Accuracy
$ws.Range("I3").Value = 0.702                   # Mean CV Test Accuracy
$ws.Range("K3").Value = "Maria"                 # Author

# --- Apply the "0.000" number format to H:I for the whole table body ----
# (rows 3-30), which is what introduces the new numFmt/cellXf in styles.xml.
$ws.Range("H3:I30").NumberFormat = "0.000"

# Re-assert the values after formatting so they keep full float precision.
$ws.Range("H3").Value = 0.709
$ws.Range("I3").Value = 0.702

# --- Column H width: 24.5 -> 26.5 (character width units) ---------------
$ws.Columns.Item(8).ColumnWidth = 26.5 - 0.8333333333333333

# --- Selection moves to G34 ----------------------------------------------
$ws.Activate()
$ws.Range("G34").Select() | Out-Null

Write-Host "Added Logistic Regression test run with scaling to Tracker sheet."
